$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column to hold the "source code" (Ma Nguon Thu) values;
# everything that used to live in columns A:F shifts right to B:G.
$ws.Columns("A:A").Insert()

# --- Row 1 (title banner) ---
# After the column insert, B1 holds the old title cell (style carried over from
# the original A1) and C1 holds the old blank companion cell. Update the title
# text in place, then move that cell back into A1 so the banner again starts
# at the left edge; B1 goes back to being completely empty.
$ws.Range("B1").Value = "THÔNG TIN IMPORT NGUỒN THU"
$ws.Range("B1").Cut($ws.Range("A1"))
$ws.Range("B1").Clear()

# --- Row 2 (column headers) ---
# New header cell for the inserted column, with its own bold/fill/border style.
$ws.Range("A2").Value = "Mã Nguồn Thu"
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Interior.Color = 49407
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("A2").Borders.ColorIndex = 64
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4108

# --- Row 3 / Row 4 (sample data rows demonstrating the import codes) ---
$ws.Range("A3").Value = "DA0001"
$ws.Range("B3").Value = "Test dự án 1"
$ws.Range("C3").Value = "điều luậ 20"
$ws.Range("D3").Value = "01/01/2025"
$ws.Range("E3").Value = "01/01/2025"
$ws.Range("F3").Value = 5000000

$ws.Range("A4").Value = "DA0002"
$ws.Range("B4").Value = "test dự án 2"
$ws.Range("C4").Value = "điều luậ 21"
$ws.Range("D4").Value = "01/01/2025"
$ws.Range("E4").Value = "01/01/2025"
$ws.Range("F4").Value = 700000

# Keep the selection sensible after these edits.
$ws.Range("B9").Select()
